$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.568.33"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "2.347.56"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.52"
$ws.Range("E5").Value = "  -4.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.66"
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  -1.99%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.616"
$ws.Range("E9").Value = "  -6.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.28"
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.63"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.996"
$ws.Range("E13").Value = "  -4.78%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.106"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.96"
$ws.Range("E15").Value = "  -7.32%  "
$ws.Range("D16").Value = "2.703.84"
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").Value = "2.354.72"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").Value = "42.560.59"
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.65"
$ws.Range("E19").Value = "  -3.79%  "
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.84"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.62"
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.82"
$ws.Range("E23").Value = "  -6.88%  "
$ws.Range("E24").Value = "  -4.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.40"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.39"
$ws.Range("E27").Value = "  -3.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.70"
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "173.03"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.86"
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0889"
$ws.Range("E32").Value = "  -4.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.05"
$ws.Range("E33").Value = "  +2.98%  "
$ws.Range("E34").Value = "  -9.69%  "
$ws.Range("E35").Value = "  +16.66%  "
$ws.Range("E36").Value = "  -2.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.62"
$ws.Range("E37").Value = "  -6.20%  "
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.95"
$ws.Range("E39").Value = "  -6.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.66"
$ws.Range("E40").Value = "  -5.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.237"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.47"
$ws.Range("E42").Value = "  -7.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.68"
$ws.Range("E43").Value = "  +1.31%  "
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.02"
$ws.Range("E45").Value = "  -3.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "111.52"
$ws.Range("E46").Value = "  -8.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.17"
$ws.Range("E47").Value = "  -2.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.45"
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.37"
$ws.Range("E49").Value = "  -6.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.93"
$ws.Range("E50").Value = "  +2.17%  "
$ws.Range("E51").Value = "  -2.27%  "
